$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap rows 12 and 13 (K_URBAN <-> K_TYPEAREA)
$ws.Range("A12").Value = "K_TYPEAREA"
$ws.Range("B12").Value = "Art der Fläche"
$ws.Range("C12").Value = "Type of area"

$ws.Range("A13").Value = "K_URBAN"
$ws.Range("B13").Value = "Verstädterungsgrad"
$ws.Range("C13").Value = "Degree of urbanisation"

# Add two new rows (14, 15), carrying over the same formatting as row 13
$ws.Range("A13:C13").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A15:C15").PasteSpecial(-4122)

$ws.Range("A14").Value = "K_SUBINDEX"
$ws.Range("B14").Value = "Subindikatoren"
$ws.Range("C14").Value = "Sub index"

$ws.Range("A15").Value = "K_CRIMOFF"
$ws.Range("B15").Value = "Straftaten"
$ws.Range("C15").Value = "Criminal offences"
